$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: date/vol/prices change (swap to what used to be row 5's data)
$ws.Range("D2").Value = 44874
$ws.Range("M2").Value = 200
$ws.Range("P2").Value = 7750
$ws.Range("S2").Value = 7750

# Row 4: date/quality/vol/prices change (swap to what used to be row 2/3's data)
$ws.Range("D4").Value = 44923
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 7500
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 7625
$ws.Range("S4").Value = 7625

# Row 5: date/quality/vol/prices change (swap to what used to be row 4's data)
$ws.Range("D5").Value = 44881
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 11250
$ws.Range("O5").Value = 11250
$ws.Range("P5").Value = 11250
$ws.Range("S5").Value = 11250
